# Apply crypto price/volume updates (and two row-pair reorderings) per the commit diff.
# Values are written as text (not auto-coerced to numbers) by temporarily forcing a
# "@" (Text) number format and restoring the original style afterwards, so cell
# styling stays byte-identical to the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$s = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.141.50'
$ws.Range("D2").Style = $s
$s = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.73%  '
$ws.Range("E2").Style = $s

# Row 3
$s = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.041.37'
$ws.Range("D3").Style = $s
$s = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.62%  '
$ws.Range("E3").Style = $s

# Row 4
$s = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E4").Style = $s

# Row 5
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.23'
$ws.Range("D5").Style = $s
$s = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E5").Style = $s

# Row 6
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.73'
$ws.Range("D6").Style = $s
$s = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +7.42%  '
$ws.Range("E6").Style = $s

# Row 7
$s = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E7").Style = $s

# Row 8
$s = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.041.06'
$ws.Range("D8").Style = $s
$s = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.70%  '
$ws.Range("E8").Style = $s

# Row 9
$s = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("E9").Style = $s

# Row 10
$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.87'
$ws.Range("D10").Style = $s
$s = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +13.37%  '
$ws.Range("E10").Style = $s

# Row 11
$s = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.42%  '
$ws.Range("E11").Style = $s

# Row 12
$s = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("E12").Style = $s

# Row 13
$s = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000235'
$ws.Range("D13").Style = $s
$s = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.81%  '
$ws.Range("E13").Style = $s

# Row 14
$s = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.86'
$ws.Range("D14").Style = $s
$s = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.26%  '
$ws.Range("E14").Style = $s

# Row 15
$s = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("E15").Style = $s

# Row 16
$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.547.46'
$ws.Range("D16").Style = $s
$s = $ws.Range("E16").Style
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("E16").Style = $s

# Row 17
$s = $ws.Range("B17").Style
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("B17").Style = $s
$s = $ws.Range("C17").Style
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C17").Style = $s
$s = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.10'
$ws.Range("D17").Style = $s
$s = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.98%  '
$ws.Range("E17").Style = $s

# Row 18
$s = $ws.Range("B18").Style
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("B18").Style = $s
$s = $ws.Range("C18").Style
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C18").Style = $s
$s = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.143.58'
$ws.Range("D18").Style = $s
$s = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.73%  '
$ws.Range("E18").Style = $s

# Row 19
$s = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.043.98'
$ws.Range("D19").Style = $s
$s = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("E19").Style = $s

# Row 20
$s = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("E20").Style = $s

# Row 21
$s = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.32'
$ws.Range("D21").Style = $s
$s = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.26%  '
$ws.Range("E21").Style = $s

# Row 22
$s = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.28%  '
$ws.Range("E22").Style = $s

# Row 23
$s = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.00%  '
$ws.Range("E23").Style = $s

# Row 24
$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.11'
$ws.Range("D24").Style = $s
$s = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.93%  '
$ws.Range("E24").Style = $s

# Row 25
$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.24'
$ws.Range("D25").Style = $s
$s = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.35%  '
$ws.Range("E25").Style = $s

# Row 26
$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").Style = $s
$s = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.29%  '
$ws.Range("E26").Style = $s

# Row 27
$s = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.38'
$ws.Range("D27").Style = $s
$s = $ws.Range("E27").Style
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.56%  '
$ws.Range("E27").Style = $s

# Row 28
$s = $ws.Range("E28").Style
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E28").Style = $s

# Row 29
$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.50'
$ws.Range("D29").Style = $s
$s = $ws.Range("E29").Style
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.19%  '
$ws.Range("E29").Style = $s

# Row 30
$s = $ws.Range("B30").Style
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("B30").Style = $s
$s = $ws.Range("C30").Style
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C30").Style = $s
$s = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.24'
$ws.Range("D30").Style = $s
$s = $ws.Range("E30").Style
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.91%  '
$ws.Range("E30").Style = $s

# Row 31
$s = $ws.Range("B31").Style
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("B31").Style = $s
$s = $ws.Range("C31").Style
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C31").Style = $s
$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.70'
$ws.Range("D31").Style = $s
$s = $ws.Range("E31").Style
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("E31").Style = $s

# Row 32
$s = $ws.Range("E32").Style
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("E32").Style = $s

# Row 33
$s = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.66'
$ws.Range("D33").Style = $s
$s = $ws.Range("E33").Style
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("E33").Style = $s

# Row 34
$s = $ws.Range("E34").Style
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("E34").Style = $s

# Row 35
$s = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0863'
$ws.Range("D35").Style = $s
$s = $ws.Range("E35").Style
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.82%  '
$ws.Range("E35").Style = $s

# Row 36
$s = $ws.Range("E36").Style
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.04%  '
$ws.Range("E36").Style = $s

# Row 37
$s = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.95'
$ws.Range("D37").Style = $s
$s = $ws.Range("E37").Style
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.01%  '
$ws.Range("E37").Style = $s

# Row 38
$s = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +12.03%  '
$ws.Range("E38").Style = $s

# Row 39
$s = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.12%  '
$ws.Range("E39").Style = $s

# Row 40
$s = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("D40").Style = $s
$s = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.08%  '
$ws.Range("E40").Style = $s

# Row 41
$s = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E41").Style = $s

# Row 42
$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.13'
$ws.Range("D42").Style = $s
$s = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("E42").Style = $s

# Row 43
$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.303'
$ws.Range("D43").Style = $s
$s = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +11.94%  '
$ws.Range("E43").Style = $s

# Row 44
$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.75'
$ws.Range("D44").Style = $s
$s = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +9.64%  '
$ws.Range("E44").Style = $s

# Row 45
$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '396.42'
$ws.Range("D45").Style = $s
$s = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("E45").Style = $s

# Row 46
$s = $ws.Range("E46").Style
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.63%  '
$ws.Range("E46").Style = $s

# Row 47
$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.729.28'
$ws.Range("D47").Style = $s
$s = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("E47").Style = $s

# Row 48
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.01'
$ws.Range("D48").Style = $s
$s = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("E48").Style = $s

# Row 49
$s = $ws.Range("B49").Style
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'USDe'
$ws.Range("B49").Style = $s
$s = $ws.Range("C49").Style
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("C49").Style = $s
$s = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = $s
$s = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E49").Style = $s

# Row 50
$s = $ws.Range("B50").Style
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("B50").Style = $s
$s = $ws.Range("C50").Style
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("C50").Style = $s
$s = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("D50").Style = $s
$s = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +6.87%  '
$ws.Range("E50").Style = $s

# Row 51
$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.55'
$ws.Range("D51").Style = $s
$s = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.32%  '
$ws.Range("E51").Style = $s
